$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "police trop petite" action recommandée changes
$ws.Range("D4").Value = "augmenter a 12px"

# Row 5: remove the "cibles tactiles trop petites / 48px" entry (B5:D5), keep A5 = SEO
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = $null

# Row 6: shift the "accessibilité / contraste" block up one row (was row6 without A,
# content used to start at A6 with old text referencing contrast, but category label
# "accessibilité" needs to appear starting row 6 instead of empty row 6/ existing row8)
$ws.Range("A6").Value = "accessibilité"
$ws.Range("B6").Value = "couleurs d'arrière-plan et de premier plan n'ont pas un rapport de contraste suffisant"
$ws.Range("C6").Value = "pas un rapport de contraste suffisamment élevé"
$ws.Range("D6").Value = "Le texte de 18 points ou 14 points en gras nécessite un rapport de contraste de 3: 1."

# Row 7: only D7 retained
$ws.Range("D7").Value = "Tout autre texte a besoin d'un rapport de contraste de 4,5: 1."

# Row 8: header ordering issue block
$ws.Range("A8").Value = "accessibilité"
$ws.Range("B8").Value = "Les éléments d'en-tête ne sont pas dans un ordre séquentiel décroissant"
$ws.Range("C8").Value = "Échec de l'audit des niveaux de titre du phare "
$ws.Range("D8").Value = "utiliser les h1 h2 h3 etc,,,"

# Row 9: just the category label
$ws.Range("A9").Value = "accessibilité"

# Clear the boolean checkbox values in column E (rows 2-15) -- keep formatting/style only
$ws.Range("E2:E15").ClearContents()

# Update the sheet view: zoom to 60%, and move selection to A10
$ws.Application.ActiveWindow.Zoom = 60
$ws.Range("A10").Select()

# Autofit columns B:E to reflect new content widths
$ws.Range("B:E").Columns.AutoFit()
